$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "27.347.61"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.79%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.708.09"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -0.96%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "224.11"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -0.77%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5305"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.68%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2663"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.75%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06614"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  +0.01%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.79"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07672"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -0.74%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.516"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -2.20%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.942.84"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.96%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.682.46"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -2.54%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5822"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -0.73%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8181"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.64"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "27.343.47"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.84%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -2.77%  "
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.631"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.21%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.41"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.64%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.991"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -1.78%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -0.15%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "143.78"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -3.02%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.686"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1203"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.75%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "7.243"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -2.32%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "16.23"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -2.44%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05374"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.30%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.485"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.82%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.429"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -0.79%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -0.76%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.861"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +1.27%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9506"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.42%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.396"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -2.05%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5863"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.58%  "
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -0.57%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.814"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -2.18%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.047.21"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.87%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8432"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -1.33%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.13%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.98"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.53%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.850.97"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -0.90%  "
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +2.72%  "
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.12%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.4522"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.78%  "
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "8.074"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -1.46%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05235"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.82%  "
